{"js": "// Remove the trailing \"Ver no Jupiter ...\" / \"\u00a9 2020 ...\" footer block\n// (and the blank paragraph right before it) that used to follow the\n// \"LOM3094: ... (Requisito)\" paragraph, while leaving the rest of the\n// document (including the blank paragraph / page-break paragraph that\n// follow the footer) untouched.\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items/text\");\nawait context.sync();\n\nconst marker = \"LOM3094: Processamento de Materiais Met\u00e1licos II (Requisito)\";\nconst footerTexts = [\n  \"Ver no Jupiter Salvar em pdf Salvar em docx\",\n  \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\",\n];\n\nconst items = paragraphs.items;\nlet markerIndex = -1;\nfor (let i = 0; i < items.length; i++) {\n  if (items[i].text === marker) {\n    markerIndex = i;\n    break;\n  }\n}\n\nif (markerIndex === -1) {\n  throw new Error(\"Could not find the requisito marker paragraph\");\n}\n\n// The paragraph immediately after the marker is the blank separator\n// paragraph, followed by the two footer paragraphs. Delete those three.\nconst toDelete = [];\nif (markerIndex + 1 < items.length && items[markerIndex + 1].text === \"\") {\n  toDelete.push(items[markerIndex + 1]);\n}\nfor (let i = markerIndex + 2; i < items.length; i++) {\n  if (footerTexts.includes(items[i].text)) {\n    toDelete.push(items[i]);\n  }\n  if (toDelete.length >= 1 + footerTexts.length) break;\n}\n\nfor (const p of toDelete) {\n  p.delete();\n}\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n\n# Locate the \"LOM3094: ... (Requisito)\" paragraph; the footer block we need\n# to remove (a blank separator paragraph, the \"Ver no Jupiter ...\" line and\n# the \"\u00a9 2020 ...\" copyright line) immediately follows it.\n$marker = \"LOM3094: Processamento de Materiais Met\u00e1licos II (Requisito)\"\n$count = $d.Paragraphs.Count\n$markerIndex = -1\nfor ($i = 1; $i -le $count; $i++) {\n    $text = $d.Paragraphs.Item($i).Range.Text.Trim()\n    if ($text -eq $marker) {\n        $markerIndex = $i\n        break\n    }\n}\n\nif ($markerIndex -eq -1) {\n    throw \"Could not find the requisito marker paragraph\"\n}\n\n$footerTexts = @(\n    \"Ver no Jupiter Salvar em pdf Salvar em docx\",\n    \"\u00a9 2020 . Contact: luizeleno@usp.br. Powered by Jekyll and Github pages. Original theme under Creative Commons Attribution\"\n)\n\n# Delete the blank paragraph right after the marker, plus the following two\n# footer paragraphs. Deleting repeatedly at ($markerIndex + 1) works because\n# each deletion shifts the subsequent paragraphs up by one.\n$toRemove = 1 + $footerTexts.Count\nfor ($n = 0; $n -lt $toRemove; $n++) {\n    $d.Paragraphs.Item($markerIndex + 1).Range.Delete()\n}\n"}
